$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to make edits, then restore protection.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclosure note (cell A38).
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for the holdings rows.
$ws.Range("D2").Value = 0.0364221309654531
$ws.Range("E2").Value = -0.00115473441108549
$ws.Range("D3").Value = 0.02044655693064605
$ws.Range("E3").Value = 0.000777907429015956
$ws.Range("D4").Value = 0.01931787615581907
$ws.Range("E4").Value = -0.003208985158443589
$ws.Range("D5").Value = 0.03770935723304032
$ws.Range("E5").Value = 0.0059753954305799
$ws.Range("D6").Value = 0.03514749657462208
$ws.Range("E6").Value = 0.002731174404994263
$ws.Range("D7").Value = 0.01985514655716964
$ws.Range("E7").Value = -0.001347968419025625
$ws.Range("D8").Value = 0.03667218422654282
$ws.Range("E8").Value = 0.003016591251885581
$ws.Range("D9").Value = 0.02045022743723085
$ws.Range("E9").Value = -0.004576864399174485
$ws.Range("D10").Value = 0.02545139461754141
$ws.Range("E10").Value = -0.006509764646970506
$ws.Range("D11").Value = 0.02344408632897718
$ws.Range("E11").Value = 0.001902690948627273
$ws.Range("D12").Value = 0.0568550764341805
$ws.Range("E12").Value = -0.001436437634666055
$ws.Range("D13").Value = 0.02506955997420564
$ws.Range("E13").Value = 0.005124450951683679
$ws.Range("D14").Value = 0.02730989445164675
$ws.Range("E14").Value = -0.01025562528700441
$ws.Range("D15").Value = 0.03271644869254556
$ws.Range("E15").Value = -0.006282722513088923
$ws.Range("D16").Value = 0.01911732375436608
$ws.Range("E16").Value = 0.0005173305742369738
$ws.Range("D17").Value = 0.03025873865847111
$ws.Range("E17").Value = 0.005896723106730617
$ws.Range("D18").Value = 0.04241525255030672
$ws.Range("E18").Value = 0.001370175839232601
$ws.Range("D19").Value = 0.1263152842316725
$ws.Range("E19").Value = 0.00265076209410231
$ws.Range("D20").Value = 0.009131608631892814
$ws.Range("E20").Value = -0.004494093477144134
$ws.Range("D21").Value = 0.01558420627020242
$ws.Range("E21").Value = -0.004671292161846585
$ws.Range("D22").Value = 0.01680255956006646
$ws.Range("E22").Value = 0.004590468330926356
$ws.Range("D23").Value = 0.01603119240541845
$ws.Range("E23").Value = 0.006585788561525252
$ws.Range("D24").Value = 0.02152650153470814
$ws.Range("E24").Value = 0.00162932790224013
$ws.Range("D25").Value = 0.01185558333113987
$ws.Range("E25").Value = 0.006914433880725879
$ws.Range("D26").Value = 0.041706131069826
$ws.Range("E26").Value = -0.01094609460946105
$ws.Range("D27").Value = 0.02390774268159919
$ws.Range("E27").Value = 0.0002942618930847907
$ws.Range("D28").Value = 0.04601867043142099
$ws.Range("E28").Value = 0.0009438414346389923
$ws.Range("D29").Value = 0.05606484491485694
$ws.Range("E29").Value = 0.003384417422175101
$ws.Range("D30").Value = 0.01292232430734828
$ws.Range("E30").Value = 0.01047806155861175
$ws.Range("D31").Value = 0.0206078553033449
$ws.Range("E31").Value = -0.0007668711656441118
$ws.Range("D32").Value = 0.01442713004860134
$ws.Range("E32").Value = 0.002650176678445249
$ws.Range("D33").Value = 0.04181303457410838
$ws.Range("E33").Value = -0.000515729757606942
$ws.Range("D34").Value = 0.01662657916102839
$ws.Range("E34").Value = 0.003004807692307709
$ws.Range("E35").Value = 0.0004006480287530589

$ws.Protect("D382")
